$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (LinearRegression) - only C2 and D2 change
$ws.Range("C2").Value = 0.5730332030979277
$ws.Range("D2").Value = 0.5730332030979277

# Row 3 (RandomForestRegressor) - values change
$ws.Range("B3").Value = 0.8450135306879822
$ws.Range("C3").Value = 0.8459775437763802
$ws.Range("D3").Value = 0.7902435973911706

# Row 4: GradientBoostingRegressor -> DecisionTreeRegressor
$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("B4").Value = 0.806352936269886
$ws.Range("C4").Value = 0.804006689235123
$ws.Range("D4").Value = 0.7825155647311234

# Row 5: AdaBoostRegressor -> MLPRegressor
$ws.Range("A5").Value = "MLPRegressor"
$ws.Range("B5").Value = 0.4965509132838069
$ws.Range("C5").Value = 0.4429196630983468
$ws.Range("D5").Value = 0.2149854027088011
